# Apply the commit's changes:
#  - Re-order / relabel the sampling-point rows (A = "<name>-<range>", B = "<name>")
#    on both "Chemical analysis" and "In-situ measurements" sheets.
#  - Row-height autofit bump (13.2) on header rows that still use the
#    worksheet default height.
#  - Active-sheet / selection swap: "Chemical analysis" becomes the
#    selected tab (was "In-situ measurements"), with new selections on
#    both sheets.

$wb = $excel.ActiveWorkbook

$chem   = $wb.Worksheets.Item("Chemical analysis")
$insitu = $wb.Worksheets.Item("In-situ measurements")

# New A/B values for rows 3-24 (identical re-sort on both sheets).
$rowData = @(
    @{ Row = 3;  A = "INF-1-4";      B = "INF" },
    @{ Row = 4;  A = "CW1MF01-1-5";  B = "CW1MF01" },
    @{ Row = 5;  A = "CW1MF02-1-5";  B = "CW1MF02" },
    @{ Row = 6;  A = "CW1MF05-1-4";  B = "CW1MF05" },
    @{ Row = 7;  A = "CW1MF06-1-4";  B = "CW1MF06" },
    @{ Row = 8;  A = "CW1MF09-1-5";  B = "CW1MF09" },
    @{ Row = 9;  A = "CW1MF10-1-5";  B = "CW1MF10" },
    @{ Row = 10; A = "CW1_EFF-1-4";  B = "CW1_EFF" },
    @{ Row = 11; A = "CW2MF01-1-4";  B = "CW2MF01" },
    @{ Row = 12; A = "CW2MF02-1-4";  B = "CW2MF02" },
    @{ Row = 13; A = "CW2MF05-1-4";  B = "CW2MF05" },
    @{ Row = 14; A = "CW2MF06-1-4";  B = "CW2MF06" },
    @{ Row = 15; A = "CW2MF09-1-4";  B = "CW2MF09" },
    @{ Row = 16; A = "CW2MF10-1-4";  B = "CW2MF10" },
    @{ Row = 17; A = "CW2_EFF-1-5";  B = "CW2_EFF" },
    @{ Row = 18; A = "CW3MF01-1-4";  B = "CW3MF01" },
    @{ Row = 19; A = "CW3MF02-1-4";  B = "CW3MF02" },
    @{ Row = 20; A = "CW3MF05-1-4";  B = "CW3MF05" },
    @{ Row = 21; A = "CW3MF06-1-4";  B = "CW3MF06" },
    @{ Row = 22; A = "CW3MF09-1-4";  B = "CW3MF09" },
    @{ Row = 23; A = "CW3MF10-1-4";  B = "CW3MF10" },
    @{ Row = 24; A = "CW3_EFF-1-4";  B = "CW3_EFF" }
)

foreach ($sheet in @($chem, $insitu)) {
    foreach ($entry in $rowData) {
        $sheet.Cells.Item($entry.Row, 1).Value() = $entry.A
        $sheet.Cells.Item($entry.Row, 2).Value() = $entry.B
    }
}

# Row-height bump: rows 1-3 on both sheets (and 4-5 on "Chemical analysis",
# which did not already carry an explicit height) get recalculated to 13.2
# (Excel's auto-fit height for the 10pt Arial content after the edit).
foreach ($r in 1..5) {
    $chem.Rows.Item($r).RowHeight = 13.2
}
foreach ($r in 1..3) {
    $insitu.Rows.Item($r).RowHeight = 13.2
}

# Selection / active-sheet swap: "Chemical analysis" is now the active tab
# with E18 selected; "In-situ measurements" reverts to its original
# selection A3:B24 and is no longer the active tab.
$insitu.Select() | Out-Null
$insitu.Range("A3:B24").Select() | Out-Null

$chem.Select() | Out-Null
$chem.Range("E18").Select() | Out-Null
